$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-07 Sunday", "2025-09-08 Monday"),
    @("38×78=2964", "33×94=3102"),
    @("68×42=2856", "79×38=3002"),
    @("13×46=598", "86×16=1376"),
    @("34×85=2890", "80×15=1200"),
    @("89×98=8722", "70×91=6370"),
    @("28×68=1904", "63×35=2205"),
    @("66×12=792", "91×21=1911"),
    @("58×51=2958", "62×94=5828"),
    @("55×49=2695", "96×67=6432"),
    @("87×99=8613", "12×96=1152"),
    @("75×21=1575", "35×79=2765"),
    @("98×50=4900", "56×98=5488"),
    @("19×42=798", "18×62=1116"),
    @("85×63=5355", "73×96=7008"),
    @("68×35=2380", "20×87=1740"),
    @("23×34=782", "94×44=4136"),
    @("48×87=4176", "70×33=2310"),
    @("23×80=1840", "56×22=1232"),
    @("25×45=1125", "43×18=774"),
    @("64×19=1216", "79×83=6557"),
    @("29×38=1102", "53×67=3551"),
    @("58×37=2146", "39×79=3081"),
    @("62×31=1922", "51×95=4845"),
    @("66×64=4224", "31×80=2480"),
    @("14×72=1008", "69×52=3588")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
